# Add the "2022-Q3" quarterly sheet (with the new holdings snapshot) right
# after the "总计" (summary) sheet, and record its summary row at the top
# of the "总计" sheet's data (below the header).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Header row (bold, centered, thin-bordered like every other quarter sheet).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # headers start at column B
    $cell = $q3.Cells.Item(1, $col)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $cell.Value = $headers[$i]
}

# Data rows: index column (A) + rank column (H) are numeric, everything
# else (fund code / name / scale / position / rank-of-position / holding
# value) is stored as text, matching the rest of the workbook.
$q3Rows = @(
    @("160135", "南方中证高铁产业指数（LOF）", "1.84", "95.01", "2.59", "0.0477", 9),
    @("160639", "鹏华中证高铁产业指数（LOF）A", "0.75", "94.62", "2.57", "0.0193", 9),
    @("015678", "鹏华中证高铁产业指数（LOF）C", "0.06", "94.62", "2.57", "0.0015", 9)
)

for ($r = 0; $r -lt $q3Rows.Length; $r++) {
    $row = 2 + $r
    $data = $q3Rows[$r]

    $idxCell = $q3.Cells.Item($row, 1)
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1
    $idxCell.Value = $r

    $q3.Range("B" + $row + ":G" + $row).NumberFormat = "@"
    $q3.Cells.Item($row, 2).Value = $data[0]
    $q3.Cells.Item($row, 3).Value = $data[1]
    $q3.Cells.Item($row, 4).Value = $data[2]
    $q3.Cells.Item($row, 5).Value = $data[3]
    $q3.Cells.Item($row, 6).Value = $data[4]
    $q3.Cells.Item($row, 7).Value = $data[5]
    $q3.Cells.Item($row, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2. Insert the matching "2022-Q3" summary row into "总计", just below
#    its header row, pushing every older quarter down by one row.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$a2 = $total.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
$a2.Value = 0

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.07000000000000001

# The row-insert kept each older row's original 0-based index value in
# column A (e.g. the old "2022-Q2" row stayed "0"); renumber column A for
# every data row so it again reflects its (0-based) position top-to-bottom.
$lastRow = $total.UsedRange.Rows.Count()
for ($r = 3; $r -le $lastRow; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
